$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("planes")

# --- 1. Update the id -> color lookup table (K1:L20). The palette of 19
#        per-id colors is collapsed down to just 3 distinct colors. ---
$ws.Range("L2").Value = "#63b8ff"
$ws.Range("L3").Value = "#8ee5ee"
$ws.Range("L4").Value = "#8ee5ee"
$ws.Range("L5").Value = "#8ee5ee"
$ws.Range("L6").Value = "#8ee5ee"
$ws.Range("L7").Value = "#8ee5ee"
$ws.Range("L8").Value = "#8ee5ee"
$ws.Range("L9").Value = "#8ee5ee"
$ws.Range("L10").Value = "#8ee5ee"
$ws.Range("L11").Value = "#8ee5ee"
$ws.Range("L12").Value = "#8ee5ee"
$ws.Range("L13").Value = "#8ee5ee"
$ws.Range("L14").Value = "#8ee5ee"
$ws.Range("L15").Value = "#8ee5ee"
$ws.Range("L16").Value = "#8ee5ee"
$ws.Range("L17").Value = "#8ee5ee"
$ws.Range("L18").Value = "#8ee5ee"
$ws.Range("L19").Value = "#dda0dd"
$ws.Range("L20").Value = "#63b8ff"

# --- 2. Replace the static color values in column B (rows 2-240) with a
#        VLOOKUP formula that reads the color from the id -> color table.
#        (Filled in chunks, matching how Excel grouped the shared formulas.) ---
$ws.Range("B2").Formula = "=VLOOKUP(A2,`$K`$1:`$L`$20,2,FALSE)"
$ws.Range("B3:B66").Formula = "=VLOOKUP(A3,`$K`$1:`$L`$20,2,FALSE)"
$ws.Range("B67:B130").Formula = "=VLOOKUP(A67,`$K`$1:`$L`$20,2,FALSE)"
$ws.Range("B131:B194").Formula = "=VLOOKUP(A131,`$K`$1:`$L`$20,2,FALSE)"
$ws.Range("B195:B240").Formula = "=VLOOKUP(A195,`$K`$1:`$L`$20,2,FALSE)"

# --- 3. Update the selected cell in the sheet view. ---
$ws.Range("N13").Select()
